# "deep sea double count fix"
# Recomputed percent-coverage / landings figures after excluding ISSCAAP
# code 46 (Deep Sea) landings from the FAO-area "double count" fix, and
# updated the footnote to mention 'Deep Sea' among the groups whose
# landings are incorporated into FAO major fishing areas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = 1.422764467196325
$ws.Range("D4").Value = 15.10956701022704
$ws.Range("E4").Value = 69.04212945962217
$ws.Range("F4").Value = 15.84830353015079
$ws.Range("G4").Value = 84.15169646984921
$ws.Range("H4").Value = 15.84830353015079

$ws.Range("B5").Value = 8.03286593
$ws.Range("C5").Value = 7.33419037196252
$ws.Range("D5").Value = 3.560091977063498
$ws.Range("E5").Value = 83.02406258426124
$ws.Range("F5").Value = 13.41584543867526
$ws.Range("G5").Value = 86.58415456132474
$ws.Range("H5").Value = 13.41584543867526

$ws.Range("C20").NumberFormat = "#,##0.000"
$ws.Range("C20").Value = 0.07802511163319266
$ws.Range("E20").Value = 55.08837186735091
$ws.Range("F20").Value = 44.91162813264909
$ws.Range("G20").Value = 55.08837186735091
$ws.Range("H20").Value = 44.91162813264909

$ws.Range("C22").Value = 0.05364998000000001
$ws.Range("D22").Value = 49.92385831271513
$ws.Range("E22").Value = 38.86139379735091
$ws.Range("F22").Value = 11.21474788993398
$ws.Range("G22").Value = 88.78525211006604
$ws.Range("H22").Value = 11.21474788993398

$ws.Range("B24").Value = 80.28050343000002
$ws.Range("C24").Value = 69.45807733661159
$ws.Range("D24").Value = 26.64723565892969
$ws.Range("E24").Value = 48.37613624522655
$ws.Range("F24").Value = 24.97662809584376
$ws.Range("G24").Value = 75.02337190415625
$ws.Range("H24").Value = 24.97662809584376

$ws.Range("I25").Value = "Note: Percent coverage in this sheet does not reflect reported percent coverage. For the reported percent coverage, `nthe landings of 'Deep Sea', 'Salmon', 'Tuna', and 'Sharks' are incorporated in the FAO major fishing areas `nfrom which their landings are reported. Thus, percent coverage calculated from this table will slightly different than reported elsewhere. `nArea landings exclude landings from ISSCAAP codes 61, 62, 63, 64, 71, 72, 73, 74, 81, 82, 83, 91, 92, 93, 94, `nexcept for stocks which have been incorporated in assessment."
